$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("run_1")
$ws.Cells.Item(2, 6).Value = 26.83536648750305
$ws.Cells.Item(3, 6).Value = 26.42936158180237
$ws.Cells.Item(4, 6).Value = 26.63850498199463
$ws.Cells.Item(5, 6).Value = 26.70202946662903
$ws.Cells.Item(6, 6).Value = 26.68752598762512
$ws.Cells.Item(7, 6).Value = 26.68219041824341
$ws.Cells.Item(8, 6).Value = 26.62762475013733
$ws.Cells.Item(9, 6).Value = 26.79902958869934
$ws.Cells.Item(10, 6).Value = 26.72500348091125
$ws.Cells.Item(11, 6).Value = 26.833420753479
$ws.Cells.Item(12, 6).Value = 26.70390462875367
$ws.Cells.Item(13, 6).Value = 26.691734790802
$ws.Cells.Item(14, 6).Value = 26.67958521842957
$ws.Cells.Item(15, 6).Value = 26.70254111289978
$ws.Cells.Item(16, 6).Value = 26.72662425041199
$ws.Cells.Item(17, 6).Value = 26.53722834587097
$ws.Cells.Item(18, 6).Value = 26.6410915851593
$ws.Cells.Item(19, 6).Value = 26.6705129146576
$ws.Cells.Item(20, 6).Value = 26.56576704978943
$ws.Cells.Item(21, 6).Value = 26.98866891860962

$ws = $wb.Worksheets.Item("run_2")
$ws.Cells.Item(2, 6).Value = 26.96417093276977
$ws.Cells.Item(3, 6).Value = 26.66024994850159
$ws.Cells.Item(4, 6).Value = 26.71920442581177
$ws.Cells.Item(5, 6).Value = 26.68315935134888
$ws.Cells.Item(6, 6).Value = 26.75304913520813
$ws.Cells.Item(7, 6).Value = 26.72948503494263
$ws.Cells.Item(8, 6).Value = 26.66810154914856
$ws.Cells.Item(9, 6).Value = 26.74563002586365
$ws.Cells.Item(10, 6).Value = 26.67763829231263
$ws.Cells.Item(11, 6).Value = 26.88343358039856
$ws.Cells.Item(12, 6).Value = 26.79745745658875
$ws.Cells.Item(13, 6).Value = 26.66997885704041
$ws.Cells.Item(14, 6).Value = 26.76531839370728
$ws.Cells.Item(15, 6).Value = 26.78905320167541
$ws.Cells.Item(16, 6).Value = 26.70238089561462
$ws.Cells.Item(17, 6).Value = 26.56805348396301
$ws.Cells.Item(18, 6).Value = 26.67107510566711
$ws.Cells.Item(19, 6).Value = 26.56411337852478
$ws.Cells.Item(20, 6).Value = 26.58590602874756
$ws.Cells.Item(21, 6).Value = 26.91187143325806

$ws = $wb.Worksheets.Item("run_3")
$ws.Cells.Item(2, 6).Value = 27.03681397438049
$ws.Cells.Item(3, 6).Value = 26.54305362701416
$ws.Cells.Item(4, 6).Value = 26.61024308204651
$ws.Cells.Item(5, 6).Value = 26.72766065597535
$ws.Cells.Item(6, 6).Value = 26.64059782028198
$ws.Cells.Item(7, 6).Value = 26.57769894599915
$ws.Cells.Item(8, 6).Value = 26.53677916526794
$ws.Cells.Item(9, 6).Value = 26.63611721992493
$ws.Cells.Item(10, 6).Value = 26.7058162689209
$ws.Cells.Item(11, 6).Value = 26.82839226722717
$ws.Cells.Item(12, 6).Value = 26.58232665061951
$ws.Cells.Item(13, 6).Value = 26.64230847358704
$ws.Cells.Item(14, 6).Value = 26.58404636383057
$ws.Cells.Item(15, 6).Value = 26.75762128829956
$ws.Cells.Item(16, 6).Value = 26.67815589904785
$ws.Cells.Item(17, 6).Value = 26.78945422172546
$ws.Cells.Item(18, 6).Value = 26.52388954162598
$ws.Cells.Item(19, 6).Value = 26.53056931495667
$ws.Cells.Item(20, 6).Value = 26.63790011405945
$ws.Cells.Item(21, 6).Value = 26.74287104606628

$ws = $wb.Worksheets.Item("run_4")
$ws.Cells.Item(2, 6).Value = 26.91141963005066
$ws.Cells.Item(3, 6).Value = 26.60600566864014
$ws.Cells.Item(4, 6).Value = 26.59395503997803
$ws.Cells.Item(5, 6).Value = 26.75175142288208
$ws.Cells.Item(6, 6).Value = 26.57131505012512
$ws.Cells.Item(7, 6).Value = 26.57181739807129
$ws.Cells.Item(8, 6).Value = 26.63414788246155
$ws.Cells.Item(9, 6).Value = 26.62747812271118
$ws.Cells.Item(10, 6).Value = 26.62723565101624
$ws.Cells.Item(11, 6).Value = 26.76723003387451
$ws.Cells.Item(12, 6).Value = 26.6086974143982
$ws.Cells.Item(13, 6).Value = 26.54006791114807
$ws.Cells.Item(14, 6).Value = 26.61694264411926
$ws.Cells.Item(15, 6).Value = 26.56692719459534
$ws.Cells.Item(16, 6).Value = 26.62402892112732
$ws.Cells.Item(17, 6).Value = 26.52960801124573
$ws.Cells.Item(18, 6).Value = 26.6010799407959
$ws.Cells.Item(19, 6).Value = 26.57134222984314
$ws.Cells.Item(20, 6).Value = 26.70846152305603
$ws.Cells.Item(21, 6).Value = 26.81432223320008

$ws = $wb.Worksheets.Item("run_5")
$ws.Cells.Item(2, 6).Value = 26.72871327400208
$ws.Cells.Item(3, 6).Value = 26.61813116073608
$ws.Cells.Item(4, 6).Value = 26.51210427284241
$ws.Cells.Item(5, 6).Value = 26.59239172935486
$ws.Cells.Item(6, 6).Value = 26.53155565261841
$ws.Cells.Item(7, 6).Value = 26.75540685653687
$ws.Cells.Item(8, 6).Value = 26.40209197998047
$ws.Cells.Item(9, 6).Value = 26.59920716285706
$ws.Cells.Item(10, 6).Value = 26.6881160736084
$ws.Cells.Item(11, 6).Value = 26.88089323043823
$ws.Cells.Item(12, 6).Value = 26.63901042938232
$ws.Cells.Item(13, 6).Value = 26.50233554840088
$ws.Cells.Item(14, 6).Value = 26.5896999835968
$ws.Cells.Item(15, 6).Value = 26.74819898605347
$ws.Cells.Item(16, 6).Value = 26.68765830993652
$ws.Cells.Item(17, 6).Value = 26.66586971282959
$ws.Cells.Item(18, 6).Value = 26.59842848777771
$ws.Cells.Item(19, 6).Value = 26.65492725372314
$ws.Cells.Item(20, 6).Value = 26.60842227935791
$ws.Cells.Item(21, 6).Value = 26.85299611091614
